# Menambah Semester & KKM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update score values (Tugas / UTS / UAS) for rows 2-4
$ws.Range("C2").Value = 44
$ws.Range("D2").Value = 66
$ws.Range("E2").Value = 55

$ws.Range("C3").Value = 55
$ws.Range("D3").Value = 88
$ws.Range("E3").Value = 66

$ws.Range("C4").Value = 66
$ws.Range("D4").Value = 88
$ws.Range("E4").Value = 99

# Move the active cell selection to B4
$ws.Range("B4").Select()
